$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so values can be updated.
$ws.Unprotect()

# Update the disclaimer date text in A10 (2021-05-12 -> 2021-05-13)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."
# Re-fit the row height so the embedded line break doesn't leave a stray custom row height behind.
$ws.Rows(10).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4785737418145579
$ws.Range("E2").Value = 0.005162827640984746

$ws.Range("D3").Value = 0.3402868171400062
$ws.Range("E3").Value = 0.01499952230820667

$ws.Range("D4").Value = 0.09496800847272185
$ws.Range("E4").Value = 0.01437788018433195

$ws.Range("D5").Value = 0.05415489726074501
$ws.Range("E5").Value = 0.003578850150080859

$ws.Range("D6").Value = 0.03201653531196896
$ws.Range("E6").Value = -0.002929115407147131

$ws.Range("D7").Value = 0.9999999999999998
$ws.Range("E7").Value = 0.009040404229858634

# Restore sheet protection to match the original workbook state.
$ws.Protect("", $false, $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
